$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $helper = $ws.Range("ZZ1")
    $helper.NumberFormat = "@"
    $helper.Value = $Text
    $helper.Copy()
    $ws.Range($Cell).PasteSpecial(-4163)
    $helper.Clear()
}

Set-TextValue "D2" "26.536.55"
Set-TextValue "E2" "  -7.68%  "
Set-TextValue "D3" "1.685.32"
Set-TextValue "E3" "  -6.54%  "
Set-TextValue "E4" "  +0.26%  "
Set-TextValue "D5" "216.50"
Set-TextValue "E5" "  -6.50%  "
Set-TextValue "D6" "1.005"
Set-TextValue "E6" "  +0.11%  "
Set-TextValue "D7" "0.4982"
Set-TextValue "E7" "  -16.22%  "
Set-TextValue "D8" "0.2613"
Set-TextValue "E8" "  -5.88%  "
Set-TextValue "D9" "21.69"
Set-TextValue "E9" "  -7.19%  "
Set-TextValue "D10" "0.06142"
Set-TextValue "E10" "  -10.17%  "
Set-TextValue "D11" "0.07289"
Set-TextValue "E11" "  -3.38%  "
Set-TextValue "D12" "1.680.55"
Set-TextValue "E12" "  -6.85%  "
Set-TextValue "D13" "4.431"
Set-TextValue "E13" "  -6.44%  "
Set-TextValue "D14" "0.5732"
Set-TextValue "E14" "  -8.49%  "
Set-TextValue "D15" "1.914.48"
Set-TextValue "E15" "  -6.56%  "
Set-TextValue "D16" "0.000008241"
Set-TextValue "E16" "  -11.15%  "
Set-TextValue "D17" "64.67"
Set-TextValue "E17" "  -14.12%  "
Set-TextValue "D18" "26.581.38"
Set-TextValue "E18" "  -7.40%  "
Set-TextValue "D19" "5.012"
Set-TextValue "E19" "  -8.22%  "
Set-TextValue "D20" "1.005"
Set-TextValue "E20" "  +0.12%  "
Set-TextValue "D21" "10.73"
Set-TextValue "E21" "  -6.38%  "
Set-TextValue "D22" "183.31"
Set-TextValue "E22" "  -13.11%  "
Set-TextValue "D23" "6.180"
Set-TextValue "E23" "  -9.93%  "
Set-TextValue "D24" "1.006"
Set-TextValue "E24" "  +0.17%  "
Set-TextValue "D25" "144.61"
Set-TextValue "E25" "  -6.33%  "
Set-TextValue "D26" "7.577"
Set-TextValue "E26" "  -3.38%  "
Set-TextValue "D27" "0.1132"
Set-TextValue "E27" "  -11.27%  "
Set-TextValue "D28" "15.33"
Set-TextValue "E28" "  -6.56%  "
Set-TextValue "D29" "1.315"
Set-TextValue "E29" "  -8.63%  "
Set-TextValue "D30" "0.05600"
Set-TextValue "E30" "  -9.85%  "
Set-TextValue "D31" "1.322"
Set-TextValue "E31" "  -6.97%  "
Set-TextValue "D32" "3.474"
Set-TextValue "E32" "  -8.01%  "
Set-TextValue "D33" "3.475"
Set-TextValue "E33" "  -7.46%  "
Set-TextValue "D34" "1.634"
Set-TextValue "E34" "  -5.01%  "
Set-TextValue "E35" "  -5.16%  "
Set-TextValue "D36" "2.371"
Set-TextValue "E36" "  -5.03%  "
Set-TextValue "D37" "0.5881"
Set-TextValue "E37" "  -8.24%  "
Set-TextValue "D38" "2.636"
Set-TextValue "E38" "  -3.10%  "
Set-TextValue "D39" "0.01591"
Set-TextValue "E39" "  -6.96%  "
Set-TextValue "D40" "1.074.81"
Set-TextValue "E40" "  -5.94%  "
Set-TextValue "D41" "5.921"
Set-TextValue "E41" "  -7.87%  "
Set-TextValue "D42" "0.8545"
Set-TextValue "E42" "  -1.46%  "
Set-TextValue "E43" "  -0.21%  "
Set-TextValue "D44" "98.61"
Set-TextValue "E44" "  -2.14%  "
Set-TextValue "D45" "1.842.54"
Set-TextValue "E45" "  -6.11%  "
Set-TextValue "D46" "56.42"
Set-TextValue "E46" "  -6.83%  "
Set-TextValue "E47" "  -6.41%  "
Set-TextValue "D48" "1.005"
Set-TextValue "E48" "  +0.01%  "
Set-TextValue "D49" "8.074"
Set-TextValue "E49" "  -3.12%  "
Set-TextValue "D50" "0.4328"
Set-TextValue "E50" "  -3.74%  "
Set-TextValue "D51" "0.05201"
Set-TextValue "E51" "  -4.77%  "
